$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.21%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.93%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.145"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.46%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.00%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.902"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.79%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.258"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.15%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.934"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.19%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9211"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.40%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1219"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.14%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1922"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.32%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09143"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.93%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03283"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.76%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09609"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.99%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001384"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.15%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005815"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.73%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.518"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.95%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.419"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.56%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3454"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.61%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.263"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.93%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1266"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.10%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2590"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.89%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04366"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.70%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.53%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004311"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.57%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.84%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02161"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.78%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05117"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.78%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007583"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.73%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1360"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.65%"
$ws.Range("B43").Value = "Dexo"
$ws.Range("C43").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008749"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.84%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001959"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.01%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008631"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006703"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.06%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.15%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001200"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.90%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.21%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.21%"
